$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $cCell = $ws.Cells.Item($row, 3)
    $dCell = $ws.Cells.Item($row, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    if ($cVal) {
        $cCell.Value2 = $cVal.Replace("github.com//", "github.com/")
    }

    if ($dVal) {
        $dCell.Value2 = $dVal.Replace("github.com//", "github.com/")
    }
}
